$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.047.78"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "1.810.14"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.95"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4411"
$ws.Range("E7").Value = "  +1.98%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3722"
$ws.Range("E8").Value = "  +1.43%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.72"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07692"
$ws.Range("E10").Value = "  +3.27%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.118"
$ws.Range("E11").Value = "  -0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.88"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("E15").Value = "  +1.88%  "
$ws.Range("D16").Value = "1.828.99"
$ws.Range("E16").Value = "  +3.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.69"
$ws.Range("E17").Value = "  +14.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001083"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06489"
$ws.Range("E19").Value = "  -1.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.003"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.47"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.263"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5376"
$ws.Range("E23").Value = "  -1.10%  "
$ws.Range("D24").Value = "28.099.57"
$ws.Range("E24").Value = "  +0.96%  "
$ws.Range("E25").Value = "  +2.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.093"
$ws.Range("E26").Value = "  -13.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.58"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.01"
$ws.Range("E28").Value = "  +2.74%  "
$ws.Range("D29").Value = "2.028.73"
$ws.Range("E29").Value = "  +2.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.323"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.17"
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.196"
$ws.Range("E32").Value = "  -7.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.838"
$ws.Range("E33").Value = "  +1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09239"
$ws.Range("E34").Value = "  +1.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.666"
$ws.Range("E35").Value = "  -8.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "13.06"
$ws.Range("E36").Value = "  +8.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02337"
$ws.Range("E37").Value = "  +2.02%  "
$ws.Range("E38").Value = "  -1.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.156"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6553"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06187"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.194"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.098"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.003"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.89"
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.384"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6066"
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.766"
$ws.Range("E48").Value = "  -1.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "126.74"
$ws.Range("E49").Value = "  +0.54%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.036"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.150"
$ws.Range("E51").Value = "  +2.33%  "
